$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("<down>",    "<haw>",      13),
    @("<is>",      "<is>",       21),
    @("<not>",     "<no>",       20),
    @("<three>",   "<three>",    11),
    @("<it>",      "<it>",       17),
    @("<each>",    "<each>",     20),
    @("<lima>",    "<tan>",      16),
    @("<foxtrot>", "<foxtrot>",  15),
    @("<a>",       "<a>",        13),
    @("<and>",     "<and>",      12),
    @("<is>",      "<is>",       16),
    @("<five>",    "<five>",     12),
    @("<november>","<november>", 18),
    @("<nine>",    "<nine>",     14),
    @("<number>",  "<november>", 17),
    @("<escape>",  "<escape>",   17),
    @("<but>",     "<with>",     11)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
